$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Overview" sheet: the row that used to describe 73da5f65-... moves to row 3
# (status "Ready for handoff"), and the row that used to describe
# e335aad8-... moves to row 2 (keeps status "Handed back: in sync with en-US").
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "e335aad8-f2ec-42f3-8d85-961302ae6254.md"
$wsOverview.Range("B2").Value = "e2e\e335aad8-f2ec-42f3-8d85-961302ae6254.md"

$wsOverview.Range("A3").Value = "73da5f65-62f9-4c82-8ded-08f52e849601.md"
$wsOverview.Range("B3").Value = "e2e\73da5f65-62f9-4c82-8ded-08f52e849601.md"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-02 08:55:58"

# Hyperlinks: the link target URLs stay anchored to their original cell, only
# the displayed text is swapped (matches the handed-back source data).
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/01f2955532aedbabb119020530d51d2bcf882310/e2e/73da5f65-62f9-4c82-8ded-08f52e849601.md", "", "", "e2e\e335aad8-f2ec-42f3-8d85-961302ae6254.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/01f2955532aedbabb119020530d51d2bcf882310/e2e/e335aad8-f2ec-42f3-8d85-961302ae6254.md", "", "", "e2e\73da5f65-62f9-4c82-8ded-08f52e849601.md") | Out-Null

# ---------------------------------------------------------------------------
# "zh-cn" sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "e335aad8-f2ec-42f3-8d85-961302ae6254.md"
$wsZhCn.Range("G2").Value = "e335aad8-f2ec-42f3-8d85-961302ae6254.eeec29ff922a3a968a64572c7ae7c91d4f7bdf9a.zh-cn.xlf"
$wsZhCn.Range("I2").Value = "e335aad8-f2ec-42f3-8d85-961302ae6254.md"
$wsZhCn.Range("J2").Value = "e335aad8-f2ec-42f3-8d85-961302ae6254.eeec29ff922a3a968a64572c7ae7c91d4f7bdf9a.zh-cn.xlf"

$wsZhCn.Range("A3").Value = "73da5f65-62f9-4c82-8ded-08f52e849601.md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("G3").Value = "73da5f65-62f9-4c82-8ded-08f52e849601.54c05d9722ec3b38c1d28be78353b58c30835350.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-02 08:55:53"
$wsZhCn.Range("I3").Value = "73da5f65-62f9-4c82-8ded-08f52e849601.md"
$wsZhCn.Range("J3").Value = "73da5f65-62f9-4c82-8ded-08f52e849601.54c05d9722ec3b38c1d28be78353b58c30835350.zh-cn.xlf"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/01f2955532aedbabb119020530d51d2bcf882310/e2e/73da5f65-62f9-4c82-8ded-08f52e849601.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e7a87adcad2a245e823067886939358366cac0b7/e2e/73da5f65-62f9-4c82-8ded-08f52e849601.md."

$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/01f2955532aedbabb119020530d51d2bcf882310/e2e/73da5f65-62f9-4c82-8ded-08f52e849601.md", "", "", "e335aad8-f2ec-42f3-8d85-961302ae6254.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/f690012378b9f151c423ab56aebd0280ec0955a9/e2e/73da5f65-62f9-4c82-8ded-08f52e849601.md", "", "", "e335aad8-f2ec-42f3-8d85-961302ae6254.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/01f2955532aedbabb119020530d51d2bcf882310/e2e/e335aad8-f2ec-42f3-8d85-961302ae6254.md", "", "", "73da5f65-62f9-4c82-8ded-08f52e849601.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/f690012378b9f151c423ab56aebd0280ec0955a9/e2e/e335aad8-f2ec-42f3-8d85-961302ae6254.md", "", "", "73da5f65-62f9-4c82-8ded-08f52e849601.md") | Out-Null

# ---------------------------------------------------------------------------
# "de-de" sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "e335aad8-f2ec-42f3-8d85-961302ae6254.md"
$wsDeDe.Range("G2").Value = "e335aad8-f2ec-42f3-8d85-961302ae6254.eeec29ff922a3a968a64572c7ae7c91d4f7bdf9a.de-de.xlf"
$wsDeDe.Range("I2").Value = "e335aad8-f2ec-42f3-8d85-961302ae6254.md"
$wsDeDe.Range("J2").Value = "e335aad8-f2ec-42f3-8d85-961302ae6254.eeec29ff922a3a968a64572c7ae7c91d4f7bdf9a.de-de.xlf"

$wsDeDe.Range("A3").Value = "73da5f65-62f9-4c82-8ded-08f52e849601.md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("G3").Value = "73da5f65-62f9-4c82-8ded-08f52e849601.54c05d9722ec3b38c1d28be78353b58c30835350.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-02 08:55:58"
$wsDeDe.Range("I3").Value = "73da5f65-62f9-4c82-8ded-08f52e849601.md"
$wsDeDe.Range("J3").Value = "73da5f65-62f9-4c82-8ded-08f52e849601.54c05d9722ec3b38c1d28be78353b58c30835350.de-de.xlf"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/01f2955532aedbabb119020530d51d2bcf882310/e2e/73da5f65-62f9-4c82-8ded-08f52e849601.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e7a87adcad2a245e823067886939358366cac0b7/e2e/73da5f65-62f9-4c82-8ded-08f52e849601.md."

$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/01f2955532aedbabb119020530d51d2bcf882310/e2e/73da5f65-62f9-4c82-8ded-08f52e849601.md", "", "", "e335aad8-f2ec-42f3-8d85-961302ae6254.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/9b1c3de57a29245953125eacc0825b9854541a4d/e2e/73da5f65-62f9-4c82-8ded-08f52e849601.md", "", "", "e335aad8-f2ec-42f3-8d85-961302ae6254.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/01f2955532aedbabb119020530d51d2bcf882310/e2e/e335aad8-f2ec-42f3-8d85-961302ae6254.md", "", "", "73da5f65-62f9-4c82-8ded-08f52e849601.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/9b1c3de57a29245953125eacc0825b9854541a4d/e2e/e335aad8-f2ec-42f3-8d85-961302ae6254.md", "", "", "73da5f65-62f9-4c82-8ded-08f52e849601.md") | Out-Null
